# Update meso_category_usage_frequency_relative worksheet:
# rows 2-52 get reordered (sorted ascending by column E "IA") and columns D
# (Analistas) and E (IA) get recomputed values. Columns A/B/C travel with
# their category (only row position changes for A/B/C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'outros', 0.04879594423320659, 0.05789422039054067, 0, 0),
    @(3, 'saúde', 0.007604562737642586, 0.005985673633598273, 0.004444444444444444, 0),
    @(4, 'sistema penitenciário', 0.0009505703422053232, 0.002158767539986263, 0.002222222222222222, 0.001118568232662192),
    @(5, 'nacionalidade e nacionalismo', 0.008555133079847909, 0.006476302619958787, 0, 0.001118568232662192),
    @(6, 'vida rural, vida no campo', 0.006653992395437262, 0.007751937984496124, 0.005555555555555556, 0.001118568232662192),
    @(7, 'amizade', 0.001267427122940431, 0.002256893337258365, 0.003333333333333334, 0.003355704697986577),
    @(8, 'pessoas com deficiência', 0.002534854245880862, 0.0008831321754489255, 0.003333333333333334, 0.003355704697986577),
    @(9, 'maternidade', 0.001584283903675539, 0.001471886959081543, 0.002222222222222222, 0.004474272930648769),
    @(10, 'raça e racismo', 0.0009505703422053232, 0.002256893337258365, 0.006666666666666667, 0.004474272930648769),
    @(11, 'sociedade e consumo', 0.01869455006337136, 0.01982141104896477, 0.006666666666666667, 0.005592841163310962),
    @(12, 'linguagem e representação', 0.02566539923954373, 0.02355019134530468, 0.005555555555555556, 0.005592841163310962),
    @(13, 'violência', 0.009822560202788339, 0.009420076538121873, 0.008888888888888889, 0.005592841163310962),
    @(14, 'recreação, lazer e entretenimento', 0.04372623574144487, 0.04533411833971151, 0.01, 0.006711409395973154),
    @(15, 'família', 0.008555133079847909, 0.01589637915808066, 0.04888888888888889, 0.006711409395973154),
    @(16, 'sonho e fantasia', 0.003485424588086185, 0.002158767539986263, 0.002222222222222222, 0.006711409395973154),
    @(17, 'desinformação, populismo e polarização', 0.0123574144486692, 0.01314885683446178, 0.007777777777777778, 0.007829977628635347),
    @(18, 'distopia, ficção científica e futuros imaginados', 0.02154626108998733, 0.01187322146992444, 0.005555555555555556, 0.007829977628635347),
    @(19, 'tecnologia, inovação e sociedade', 0.01045627376425855, 0.01118634088901972, 0.01444444444444444, 0.007829977628635347),
    @(20, 'questão agrária e territorial', 0.01045627376425855, 0.01246197625355706, 0.01777777777777778, 0.008948545861297539),
    @(21, 'moradia e habitação', 0.001901140684410646, 0.001766264350897851, 0.02333333333333333, 0.008948545861297539),
    @(22, 'mobilidade', 0.009822560202788339, 0.01000883132175449, 0.02333333333333333, 0.01118568232662192),
    @(23, 'territorialidade e colonialismo', 0.01077313054499366, 0.008438818565400843, 0.01777777777777778, 0.01230425055928412),
    @(24, 'violências e preconceitos de gênero', 0.005069708491761723, 0.006868805809047199, 0.004444444444444444, 0.01342281879194631),
    @(25, 'alimentação e tratamentos tradicionais', 0.009822560202788339, 0.01138259248356393, 0.008888888888888889, 0.0145413870246085),
    @(26, 'sonoridade e paisagem sonora', 0.01552598225602028, 0.0137376116180944, 0.01222222222222222, 0.01565995525727069),
    @(27, 'gênero e sexualidade', 0.02217997465145754, 0.03159650672161711, 0.01111111111111111, 0.01565995525727069),
    @(28, 'biografia', 0.003168567807351077, 0.002453144931802571, 0, 0.01677852348993289),
    @(29, 'trabalho e ofício', 0.02724968314321927, 0.02806397801982141, 0.03222222222222222, 0.01789709172259508),
    @(30, 'direitos humanos', 0.008238276299112801, 0.00863507015994505, 0.01, 0.01901565995525727),
    @(31, 'crises e desastres ambientais e sociais', 0.009188846641318124, 0.008438818565400843, 0.01555555555555556, 0.01901565995525727),
    @(32, 'religião, espiritualidade e cosmologias', 0.01204055766793409, 0.01216759886174075, 0.01555555555555556, 0.02125279642058166),
    @(33, 'solidão', 0.00285171102661597, 0.002649396526346777, 0.01666666666666667, 0.02125279642058166),
    @(34, 'saudade, luto ou perda', 0.01615969581749049, 0.01079383769993131, 0.007777777777777778, 0.02237136465324385),
    @(35, 'política', 0.02724968314321927, 0.03728780296339908, 0.03222222222222222, 0.02572706935123042),
    @(36, 'dinâmica urbana', 0.0532319391634981, 0.04189971543518791, 0.05333333333333334, 0.02572706935123042),
    @(37, 'vida afetiva', 0.09917617237008872, 0.09822392306937494, 0.06555555555555556, 0.02796420581655481),
    @(38, 'estruturas sociais e econômicas', 0.05038022813688213, 0.05102541458149348, 0.02555555555555556, 0.029082774049217),
    @(39, 'educação e socialização', 0.03295310519645121, 0.03061524874889609, 0.02888888888888889, 0.029082774049217),
    @(40, 'cultural', 0.06210392902408111, 0.06476302619958788, 0.04777777777777778, 0.0302013422818792),
    @(41, 'corpo, performance e expressão', 0.01045627376425855, 0.01530762437444804, 0.006666666666666667, 0.0302013422818792),
    @(42, 'memória e patrimônio', 0.03612167300380228, 0.03718967716612698, 0.02777777777777778, 0.0302013422818792),
    @(43, 'povos originários e comunidades tradicionais', 0.02408111533586819, 0.02237268177803945, 0.02111111111111111, 0.03243847874720358),
    @(44, 'pandemia', 0.03580481622306717, 0.02983024237071926, 0.05333333333333334, 0.03243847874720358),
    @(45, 'mulher, feminino e feminismo', 0.006653992395437262, 0.005985673633598273, 0.02444444444444445, 0.03467561521252797),
    @(46, 'vida cotidiana', 0.007604562737642586, 0.010499460308115, 0.02, 0.04138702460850112),
    @(47, 'resistência e luta', 0.00697084917617237, 0.007850063781768227, 0.03, 0.04138702460850112),
    @(48, 'saúde mental', 0.03358681875792142, 0.0278677264252772, 0.04111111111111111, 0.0436241610738255),
    @(49, 'ambiental', 0.03612167300380228, 0.03395152585614758, 0.05777777777777778, 0.05145413870246085),
    @(50, 'reflexão', 0.00285171102661597, 0.004023157688156216, 0.02111111111111111, 0.05257270693512305),
    @(51, 'poesia e ensaio', 0.07065906210392903, 0.06829555490138357, 0.06222222222222222, 0.05592841163310962),
    @(52, 'arte', 0.006337135614702154, 0.004023157688156216, 0.02666666666666667, 0.06823266219239374)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
